$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 49, shifting existing rows 49..99 down to 50..100.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with its data (same shape as the row
# that was previously there, now at row 50, with a few values updated).
$ws.Range("A49").Value = 1
$ws.Range("B49").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C49").Value = "Arica y Parinacota"
$ws.Range("D49").Value = 45210
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = 100112031
$ws.Range("G49").Value = "Poroto verde"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 1150
$ws.Range("K49").Value = 800
$ws.Range("L49").Value = 900
$ws.Range("M49").Value = 848
$ws.Range("N49").Value = "$/kilo"
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 848
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
